# "bisschen am statusbericht 4 gschaffelet"
# Fill in additional weekly hour-tracking entries (columns V/W/X/Y = weeks 18-21)
# on the "Kosten" sheet for several team members, across the
# 1_PM, 4_Realisierung, 5_Validierung and 6_Praesentationen blocks.
# All subtotal/cumulative rows (22-25, 73-76, 91-94, 109-112) are formulas
# and recalculate automatically once the raw inputs below are written.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Kosten")

# --- 1_PM block (Marina Taborda / Michel Alt), row 19-20 ---
$ws.Range("V19").Value = 2
$ws.Range("W19").Value = 3
$ws.Range("X19").Value = 3

$ws.Range("V20").Value = 1
$ws.Range("W20").Value = 1
$ws.Range("X20").Value = 1

# --- 4_Realisierung "Ist" block, rows 67-72 ---
# Marina Taborda (67): only X67 newly filled in (V67 already had a value)
$ws.Range("X67").Value = 2

# Frank Imhof (68): X/Y updated
$ws.Range("X68").Value = 15
$ws.Range("Y68").Value = 5

# Luca Krummenacher (69)
$ws.Range("V69").Value = 8
$ws.Range("W69").Value = 5
$ws.Range("X69").Value = 5

# Michel Alt (70)
$ws.Range("V70").Value = 5
$ws.Range("W70").Value = 2
$ws.Range("X70").Value = 2

# Richard Britt (71)
$ws.Range("V71").Value = 5
$ws.Range("W71").Value = 8
$ws.Range("X71").Value = 5

# Fady Angly (72)
$ws.Range("W72").Value = 1
$ws.Range("X72").Value = 1

# --- 5_Validierung "Ist" block, rows 85-89 ---
# Marina Taborda (85)
$ws.Range("V85").Value = 3
$ws.Range("W85").Value = 3
$ws.Range("X85").Value = 3
$ws.Range("Y85").Value = 3

# Luca Krummenacher (87)
$ws.Range("V87").Value = 5
$ws.Range("W87").Value = 5
$ws.Range("X87").Value = 5

# Michel Alt (88)
$ws.Range("V88").Value = 2
$ws.Range("W88").Value = 2

# Richard Britt (89)
$ws.Range("V89").Value = 5
$ws.Range("W89").Value = 1
$ws.Range("X89").Value = 1

# --- 6_Praesentationen "Ist" block, rows 103-105 ---
# Marina Taborda (103)
$ws.Range("X103").Value = 2
$ws.Range("Y103").Value = 4

# Frank Imhof (104)
$ws.Range("X104").Value = 2
$ws.Range("Y104").Value = 4

# Luca Krummenacher (105)
$ws.Range("X105").Value = 2
$ws.Range("Y105").Value = 4

# The "Status-4" sheet's view was scrolled back up near the top
# (set this first, then return to "Kosten" so it stays the active tab).
$ws4 = $wb.Worksheets.Item("Status-4")
$ws4.Range("A2").Select()

# Re-activate "Kosten" (it must remain the selected/visible tab).
$ws.Activate()

# Scroll the frozen pane down so row 101 is the first visible scrolling row,
# then restore the on-screen selection to match where work left off.
$ws.Application.Goto($ws.Range("E101"), $false)
$ws.Range("C43").Select()
